# Cards.xlsx update: "Added some new cards, also added Game Design Doc"
#
# - Card #17 (row 19) renamed from "The Plaguedoctor" to "Mad Scientist"
# - Card #20 (row 22) "Black Market" effect text split: Black Market's own
#   effect is simplified (no longer lets you kill a creature), and a couple
#   of missing HP/Strength cells (0/0) are filled in
# - Four new cards added in rows 23-26 (#21-#24): Pustulent Zombie,
#   Combat Medic, Immunity, Antidote
# - Selection moved from J27 to E27

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Card #20: Black Market effect text + fill in missing 0/0 cells ------
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = "Discard a card in your hand and gain gold equal to its cost."

# --- Card #17: The Plaguedoctor -> Mad Scientist -------------------------
$ws.Range("E19").Value = "Mad Scientist"

# --- New card #21: Pustulent Zombie (row 23) ------------------------------
$ws.Range("E23").Value = "Pustulent Zombie"
$ws.Range("F23").Value = "MINION"
$ws.Range("G23").Value = 6
$ws.Range("H23").Value = 10
$ws.Range("I23").Value = 1
$ws.Range("J23").Value = "When this minion is killed, you can bury it instead. If this minion battles an opponents minion, give it 1 poison counter."

# --- New card #22: Combat Medic (row 24) ----------------------------------
$ws.Range("E24").Value = "Combat Medic"
$ws.Range("F24").Value = "MINION"
$ws.Range("G24").Value = 2
$ws.Range("H24").Value = 2
$ws.Range("I24").Value = 2
$ws.Range("J24").Value = "When this minion is summoned, you can restore 2 health to a minion on the field. You can pay 2 gold, restore 2 health to a minion on the field."

# --- New card #23: Immunity (row 25) --------------------------------------
$ws.Range("E25").Value = "Immunity"
$ws.Range("F25").Value = "UTILITY"
$ws.Range("G25").Value = 5
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = "Remove all poison counters from all minions on your side of the field."

# --- New card #24: Antidote (row 26) --------------------------------------
$ws.Range("E26").Value = "Antidote"
$ws.Range("F26").Value = "UTILITY"
$ws.Range("G26").Value = 2
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = "Remove all poison counters from a minion on your side of the field."

# --- Move active selection from J27 to E27 --------------------------------
[void]$ws.Range("E27").Select()
